# Applies the edit described by the commit "updated task used in testing":
#   - D2: 5 -> 4
#   - F2: -2 -> -3
#   - H2: 36 -> 46
#   - active selection moves from D5 to D2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

# Move/record the active selection to D2 (was D5 before the edit)
$ws.Range("D2").Select()
